# Scheduled-runner update: refresh computed profit-margin figures
# (columns H-N) across the per-class profit sheets. Values below were
# recomputed upstream; one zero-margin outlier row (BSM!7) drops its
# now-undefined margin cell entirely.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 3752.3157
$ws.Range("I28").Value = 983.1
$ws.Range("K28").Value = 983.1
$ws.Range("M28").Value = -498.1
$ws.Range("H103").Value = 1252.8846
$ws.Range("J103").Value = 561.9524
$ws.Range("L103").Value = 1685.8572
$ws.Range("N103").Value = -2857.8572
$ws.Range("H137").Value = 5688125
$ws.Range("I137").Value = 9262151
$ws.Range("J137").Value = 11730.353
$ws.Range("K137").Value = 27786453
$ws.Range("L137").Value = 35191.05899999999
$ws.Range("M137").Value = -27783903
$ws.Range("N137").Value = -40291.05899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 52996
$ws.Range("J37").Value = 99988
$ws.Range("L37").Value = 99988
$ws.Range("N37").Value = -100534
$ws.Range("H45").Value = 1749.5
$ws.Range("I45").Value = 1749
$ws.Range("K45").Value = 1749
$ws.Range("M45").Value = -1372
$ws.Range("H132").Value = 5993.5264
$ws.Range("I132").Value = 3080.0454
$ws.Range("J132").Value = 9999.5625
$ws.Range("K132").Value = 9240.136200000001
$ws.Range("L132").Value = 29998.6875
$ws.Range("M132").Value = -6710.136200000001
$ws.Range("N132").Value = -35058.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H80").Value = 56104.555
$ws.Range("I80").Value = 774.8570999999999
$ws.Range("K80").Value = 774.8570999999999
$ws.Range("M80").Value = 223.1429000000001
$ws.Range("H83").Value = 56104.555
$ws.Range("I83").Value = 774.8570999999999
$ws.Range("K83").Value = 3874.2855
$ws.Range("M83").Value = 1117.7145
$ws.Range("H94").Value = 555.44116
$ws.Range("I94").Value = 445.17242
$ws.Range("J94").Value = 1195
$ws.Range("K94").Value = 445.17242
$ws.Range("L94").Value = 1195
$ws.Range("M94").Value = 5.827580000000012
$ws.Range("N94").Value = -2097
$ws.Range("H105").Value = 62502230
$ws.Range("I105").Value = 76925230
$ws.Range("J105").Value = 2581
$ws.Range("K105").Value = 76925230
$ws.Range("L105").Value = 2581
$ws.Range("M105").Value = -76923483
$ws.Range("N105").Value = -6075
$ws.Range("H107").Value = 1556.8636
$ws.Range("H134").Value = 4274.3184
$ws.Range("I134").Value = 2046.2667
$ws.Range("J134").Value = 9048.714
$ws.Range("K134").Value = 6138.800099999999
$ws.Range("L134").Value = 27146.142
$ws.Range("M134").Value = -3603.800099999999
$ws.Range("N134").Value = -32216.142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 305
$ws.Range("I22").Value = 360.16666
$ws.Range("J22").Value = 249.83333
$ws.Range("K22").Value = 360.16666
$ws.Range("L22").Value = 249.83333
$ws.Range("M22").Value = -10.16665999999998
$ws.Range("N22").Value = -949.8333299999999
$ws.Range("H31").Value = 1926250.5
$ws.Range("I31").Value = 1543.6333
$ws.Range("J31").Value = 4550851
$ws.Range("K31").Value = 1543.6333
$ws.Range("L31").Value = 4550851
$ws.Range("M31").Value = -1248.6333
$ws.Range("N31").Value = -4551441
$ws.Range("H34").Value = 1926250.5
$ws.Range("I34").Value = 1543.6333
$ws.Range("J34").Value = 4550851
$ws.Range("K34").Value = 1543.6333
$ws.Range("L34").Value = 4550851
$ws.Range("M34").Value = -1341.6333
$ws.Range("N34").Value = -4551255
$ws.Range("H105").Value = 10872.6
$ws.Range("I105").Value = 13215.75
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 13215.75
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = -11468.75
$ws.Range("N105").Value = -4994
$ws.Range("H107").Value = 1340.3
$ws.Range("I107").Value = 911.8333
$ws.Range("K107").Value = 911.8333
$ws.Range("M107").Value = 1008.1667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7415687.5
$ws.Range("I4").Value = 11366924
$ws.Range("J4").Value = 7117.75
$ws.Range("K4").Value = 34100772
$ws.Range("L4").Value = 21353.25
$ws.Range("M4").Value = -34100660
$ws.Range("N4").Value = -21577.25
$ws.Range("H5").Value = 322.03705
$ws.Range("J5").Value = 384.88235
$ws.Range("L5").Value = 1154.64705
$ws.Range("N5").Value = -1378.64705
$ws.Range("H76").Value = 3666.3333
$ws.Range("I76").Value = 499.5
$ws.Range("K76").Value = 1498.5
$ws.Range("M76").Value = -1115.5
$ws.Range("H79").Value = 3666.3333
$ws.Range("I79").Value = 499.5
$ws.Range("K79").Value = 1498.5
$ws.Range("M79").Value = -172.5
$ws.Range("H135").Value = 322.03705
$ws.Range("J135").Value = 384.88235
$ws.Range("L135").Value = 3463.94115
$ws.Range("N135").Value = -8533.941149999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7695.7837
$ws.Range("I132").Value = 6267.8096
$ws.Range("J132").Value = 9570
$ws.Range("K132").Value = 18803.4288
$ws.Range("L132").Value = 28710
$ws.Range("M132").Value = -16273.4288
$ws.Range("N132").Value = -33770
$ws.Range("H139").Value = 95899.60000000001
$ws.Range("J139").Value = 95899.60000000001
$ws.Range("L139").Value = 95899.60000000001
$ws.Range("N139").Value = -106179.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2633.3333
$ws.Range("H85").Value = 2633.3333
$ws.Range("H106").Value = 9995
$ws.Range("J106").Value = 9995
$ws.Range("L106").Value = 9995
$ws.Range("N106").Value = -12519
$ws.Range("H132").Value = 5870.4546
$ws.Range("I132").Value = 5643.75
$ws.Range("K132").Value = 16931.25
$ws.Range("M132").Value = -14401.25
$ws.Range("H133").Value = 50988.8
$ws.Range("J133").Value = 50988.8
$ws.Range("L133").Value = 50988.8
$ws.Range("N133").Value = -56048.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 4397.6665
$ws.Range("I61").Value = 3461.8572
$ws.Range("J61").Value = 17499
$ws.Range("K61").Value = 3461.8572
$ws.Range("L61").Value = 17499
$ws.Range("M61").Value = -3169.8572
$ws.Range("N61").Value = -18083
$ws.Range("H105").Value = 63097.8
$ws.Range("J105").Value = 63097.8
$ws.Range("L105").Value = 63097.8
$ws.Range("N105").Value = -70085.8
$ws.Range("H132").Value = 5642.6523
$ws.Range("I132").Value = 3888.1904
$ws.Range("J132").Value = 7116.4
$ws.Range("K132").Value = 11664.5712
$ws.Range("L132").Value = 21349.2
$ws.Range("M132").Value = -9134.5712
$ws.Range("N132").Value = -26409.2
$ws.Range("H138").Value = 99999
$ws.Range("J138").Value = 99999
$ws.Range("L138").Value = 99999
$ws.Range("N138").Value = -110279
